$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute(
    "Waktu Kampanye rasi bintang Taurus 2022: 16-25 Januari",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "Waktu Kampanye 2022 untuk rasi bintang Taurus: 16-25 Januari",
    2
)
